{"js": "// Office.js (Word JavaScript API) script.\n// Applies the diff: the single-column table's first three rows get their\n// value changed to \"0M\", ten new rows are inserted right after the (now)\n// third row, the two tab-separated \"summary\" rows near the end collapse\n// down to a single plain value each, and the final (previously empty) row\n// gets the text \"63\".\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Helper: replace the *entire* text content of a table cell while keeping\n// the formatting (rFonts/sz) of the first run already present in it. Using\n// body.getRange().insertText(text, \"Replace\") collapses any existing runs\n// (including ones separated by <w:tab/>) into a single run carrying the\n// original run's properties.\nfunction setCellText(cellIndex, rowIndex, text) {\n  const cell = rows.items[rowIndex].cells.items[cellIndex];\n  const range = cell.body.getRange();\n  range.insertText(text, \"Replace\");\n}\n\n// --- Step 1: change the first three rows' values to \"0M\" ---\nrows.items[0].cells.load(\"items\");\nrows.items[1].cells.load(\"items\");\nrows.items[2].cells.load(\"items\");\nawait context.sync();\n\nsetCellText(0, 0, \"0M\");\nsetCellText(0, 1, \"0M\");\nsetCellText(0, 2, \"0M\");\nawait context.sync();\n\n// --- Step 2: insert 10 new rows right after row index 2 ---\nconst newRowValues = [\n  [\"11\"],\n  [\"0.00002\"],\n  [\"0.00006\"],\n  [\"0.00003\"],\n  [\"0.00001\"],\n  [\"0.00003\"],\n  [\"0.00004\"],\n  [\"0.00004\"],\n  [\"0.00040\"],\n  [\"100.0\"],\n];\nrows.items[2].insertRows(\"After\", newRowValues.length, newRowValues);\nawait context.sync();\n\n// --- Step 3: reload rows (indices shifted by the insertion above) and fix\n// the two tab-separated \"summary\" rows + the trailing empty row. ---\nrows.load(\"items\");\nawait context.sync();\n\nconst total = rows.items.length; // 46 after the insertion\nconst idxHundred = total - 3; // was row 33 before insertion -> now +10\nconst idxOne = total - 2; // was row 34 before insertion -> now +10\nconst idxLast = total - 1; // was row 35 (empty run) before insertion -> now +10\n\nrows.items[idxHundred].cells.load(\"items\");\nrows.items[idxOne].cells.load(\"items\");\nrows.items[idxLast].cells.load(\"items\");\nawait context.sync();\n\nsetCellText(0, idxHundred, \"100\");\nsetCellText(0, idxOne, \"0\");\nsetCellText(0, idxLast, \"63\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the diff: the single-column table's first three rows get their\n# value changed to \"0M\", ten new rows are inserted right after the (now)\n# third row, the two tab-separated \"summary\" rows near the end collapse\n# down to a single plain value each, and the final (previously empty) row\n# gets the text \"63\".\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# --- Step 1: change the first three rows' values to \"0M\" ---\n# Assigning directly to Cell(...).Range.Text replaces the run's text while\n# keeping the run's existing formatting (rFonts/sz), and collapses any\n# multi-run/tab content in the cell down to a single run.\n$t.Cell(1, 1).Range.Text = \"0M\"\n$t.Cell(2, 1).Range.Text = \"0M\"\n$t.Cell(3, 1).Range.Text = \"0M\"\n\n# --- Step 2: insert 10 new rows right after row 3 ---\n# Rows.Add(beforeRow) inserts a new row immediately before beforeRow,\n# copying formatting from the neighboring row. Since row 4 (\"0\") is\n# always pushed down as we insert, repeatedly inserting \"before row 4\"\n# and filling values in reverse order lands them in the correct final\n# order (11, 0.00002, 0.00006, 0.00003, 0.00001, 0.00003, 0.00004,\n# 0.00004, 0.00040, 100.0).\n$newRowValues = @(\n    \"11\",\n    \"0.00002\",\n    \"0.00006\",\n    \"0.00003\",\n    \"0.00001\",\n    \"0.00003\",\n    \"0.00004\",\n    \"0.00004\",\n    \"0.00040\",\n    \"100.0\"\n)\nfor ($i = $newRowValues.Length - 1; $i -ge 0; $i--) {\n    $newRow = $t.Rows.Add($t.Rows.Item(4))\n    $newRow.Cells.Item(1).Range.Text = $newRowValues[$i]\n}\n\n# --- Step 3: fix the trailing three rows (indices shifted by +10 now) ---\n$total = $t.Rows.Count\n$t.Cell($total - 2, 1).Range.Text = \"100\"\n$t.Cell($total - 1, 1).Range.Text = \"0\"\n$t.Cell($total, 1).Range.Text = \"63\"\n"}
